$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 15.72
$ws.Range("E2").Value = 64.7
$ws.Range("F2").Value = 5.08
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 59.4
$ws.Range("N2").Value = 54.02451352198364

# Row 3
$ws.Range("D3").Value = 92038.3
$ws.Range("E3").Value = 61.2
$ws.Range("F3").Value = 1.31
$ws.Range("K3").Value = 54.2
$ws.Range("N3").Value = 54.02451352198364

# Row 4
$ws.Range("D4").Value = 275.46
$ws.Range("E4").Value = 46.7
$ws.Range("F4").Value = 3.96
$ws.Range("K4").Value = 51.2
$ws.Range("N4").Value = 54.02451352198364

# Row 5
$ws.Range("D5").Value = 12.41
$ws.Range("E5").Value = 47.2
$ws.Range("F5").Value = 11.7
$ws.Range("K5").Value = 49.4
$ws.Range("N5").Value = 54.02451352198364

# Row 6
$ws.Range("D6").Value = 186.65
$ws.Range("E6").Value = 39.6
$ws.Range("F6").Value = 6.27
$ws.Range("K6").Value = 35.6
$ws.Range("N6").Value = 54.02451352198364
